$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix labels: the "Z=0" / "P(Z=0)" conditions should actually be "Z=1" / "P(Z=1)"
$ws.Range("C16").Value = "P(Z=1)"
$ws.Range("C21").Value = "P(Z=1)"
$ws.Range("C11").Value = "Z=1"

# Correct the conditional probability values so that P(Y|X,Z) + P(Y|X,!Z) add up to 1 correctly
$ws.Range("D12").Value = 0.81
$ws.Range("D13").Value = 0.48
$ws.Range("D14").Value = 0.83
$ws.Range("D15").Value = 0.55

$ws.Range("D17").Value = 0.76
$ws.Range("D19").Value = 0.87
$ws.Range("D20").Value = 0.54

$ws.Range("D22").Value = 0.91
$ws.Range("D23").Value = 0.55
$ws.Range("D24").Value = 0.77
$ws.Range("D25").Value = 0.54

# Update the view to reflect where the user ended up looking
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D25").Select()
